# Add the two "Extra" sheets for Afghanistan ODI player performance data:
#   - "ODI Batting Extra"  (after "ODI Bowling")
#   - "ODI Bowling Extra"  (after "ODI Batting Extra")
#
# Header rows reuse the same bold/bordered/centered look already used by the
# header row on "Player Info" (cellXfs style index 1) by copy/pasting the
# format from an existing header cell instead of inventing a new style.

$wb = $excel.ActiveWorkbook

# Sheet used as the source of the existing bold header style.
$styleSource = $wb.Worksheets.Item("Player Info").Range("A1")

# ---------------------------------------------------------------------
# Sheet 4: "ODI Batting Extra"
# ---------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$wsBatting = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$wsBatting.Name = "ODI Batting Extra"

# Match the page-margin defaults used by the rest of the workbook
# (0.75in / 1in / 0.5in, expressed in points here).
$wsBatting.PageSetup.LeftMargin = 54
$wsBatting.PageSetup.RightMargin = 54
$wsBatting.PageSetup.TopMargin = 72
$wsBatting.PageSetup.BottomMargin = 72
$wsBatting.PageSetup.HeaderMargin = 36
$wsBatting.PageSetup.FooterMargin = 36

$wsBatting.Range("A1").Value = "MATCH_CODE"
$wsBatting.Range("B1").Value = "BATTING_POSITION"
$wsBatting.Range("C1").Value = "NUM_4"
$wsBatting.Range("D1").Value = "NUM_6"
$wsBatting.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$wsBatting.Range("F1").Value = "MAN_OF_MATCH"

$styleSource.Copy()
$wsBatting.Range("A1:F1").PasteSpecial(-4122)

# Row 2 data -- only A2 (MATCH_CODE) and F2 (MAN_OF_MATCH) are populated.
# Use a leading apostrophe so the purely-numeric "4675" is stored as text
# (matching the source data) and then strip the resulting quote-prefix
# style back to the workbook default.
$wsBatting.Cells.Item(2, 1).Formula = "'4675"
$wsBatting.Cells.Item(2, 1).Style = "Normal"
$wsBatting.Cells.Item(2, 6).Value = "NO"

# ---------------------------------------------------------------------
# Sheet 5: "ODI Bowling Extra"
# ---------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$wsBowling = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$wsBowling.Name = "ODI Bowling Extra"

# Match the page-margin defaults used by the rest of the workbook
# (0.75in / 1in / 0.5in, expressed in points here).
$wsBowling.PageSetup.LeftMargin = 54
$wsBowling.PageSetup.RightMargin = 54
$wsBowling.PageSetup.TopMargin = 72
$wsBowling.PageSetup.BottomMargin = 72
$wsBowling.PageSetup.HeaderMargin = 36
$wsBowling.PageSetup.FooterMargin = 36

$wsBowling.Range("A1").Value = "MATCH_CODE"
$wsBowling.Range("B1").Value = "MAIDEN_OVERS"
$wsBowling.Range("C1").Value = "PERCENT_WICKETS_OF_ALL"

$styleSource.Copy()
$wsBowling.Range("A1:C1").PasteSpecial(-4122)

# Row 2 data -- MATCH_CODE is populated, the other two columns are present
# but blank (empty text cells, not missing cells).
$wsBowling.Cells.Item(2, 1).Formula = "'4675"
$wsBowling.Cells.Item(2, 1).Style = "Normal"
$wsBowling.Cells.Item(2, 2).Formula = "'"
$wsBowling.Cells.Item(2, 2).Style = "Normal"
$wsBowling.Cells.Item(2, 3).Formula = "'"
$wsBowling.Cells.Item(2, 3).Style = "Normal"

# Restore the original active sheet (the workbook opened on "Player Info";
# adding sheets would otherwise leave the newly-added sheet active).
$wb.Worksheets.Item("Player Info").Activate()
